# Update total_risk (R) and total_risk_resp (S) values for rows 2-15
# per the 5 and 10 mi radius updates to the transition rule facility demographics script.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = @(38.5714285714286, 0.485714285714286)
    3  = @(50, 0.6)
    4  = @(38.2894736842105, 0.444736842105263)
    5  = @(30, 0.497674418604651)
    6  = @(30, 0.393333333333333)
    7  = @(27, 0.32)
    8  = @(20, 0.2)
    9  = @(87.5, 0.5)
    10 = @(128.75, 0.5375)
    11 = @(40, 0.5)
    12 = @(30, 0.36)
    13 = @(41.1764705882353, 0.4)
    14 = @(20, 0.21)
    15 = @(30, 0.4)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Range("R$row").Value = $pair[0]
    $ws.Range("S$row").Value = $pair[1]
}
